$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D and E hold numeric-looking / percentage-looking text that must
# remain plain text (the sheet stores them as inline strings). We briefly switch the
# cell to Text number format while assigning the value, then restore the default
# "Normal" style so the cell keeps no explicit style (matching the rest of the sheet).
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '302.33'
Set-TextValue 'E2' '-5.71%'
Set-TextValue 'D3' '35.11'
Set-TextValue 'E3' '-2.73%'
Set-TextValue 'D4' '5.048'
Set-TextValue 'E4' '-1.54%'
Set-TextValue 'D5' '0.07906'
Set-TextValue 'E5' '-3.03%'
Set-TextValue 'D6' '1.941'
Set-TextValue 'E6' '-9.72%'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue 'D7' '4.020'
Set-TextValue 'E7' '-2.89%'
$ws.Range('B8').Value = 'KuCoinToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue 'D8' '7.721'
Set-TextValue 'E8' '-3.98%'
Set-TextValue 'E9' '2.72%'
Set-TextValue 'D10' '0.9238'
Set-TextValue 'E10' '-0.34%'
Set-TextValue 'D11' '0.1197'
Set-TextValue 'E11' '18.82%'
Set-TextValue 'D12' '0.1844'
Set-TextValue 'E12' '-2.31%'
Set-TextValue 'D13' '0.09383'
Set-TextValue 'E13' '2.27%'
Set-TextValue 'D14' '0.03537'
Set-TextValue 'E14' '-1.52%'
Set-TextValue 'D15' '0.09890'
Set-TextValue 'E15' '-0.30%'
Set-TextValue 'D16' '0.001387'
Set-TextValue 'E16' '-3.12%'
Set-TextValue 'D17' '0.005821'
Set-TextValue 'E17' '2.24%'
Set-TextValue 'D18' '3.491'
Set-TextValue 'E18' '1.13%'
Set-TextValue 'E19' '2.10%'
Set-TextValue 'E20' '-0.18%'
Set-TextValue 'D21' '5.037'
Set-TextValue 'E21' '-0.42%'
Set-TextValue 'E22' '9.73%'
Set-TextValue 'D23' '0.04487'
Set-TextValue 'E23' '-2.42%'
Set-TextValue 'D24' '0.001216'
Set-TextValue 'E24' '-2.17%'
Set-TextValue 'D25' '0.004571'
Set-TextValue 'E25' '-3.34%'
Set-TextValue 'D26' '0.0001250'
Set-TextValue 'E26' '-3.83%'
Set-TextValue 'E27' '-6.82%'
Set-TextValue 'D39' '0.01906'
Set-TextValue 'E39' '-6.18%'
Set-TextValue 'E40' '-6.05%'
Set-TextValue 'D41' '0.007597'
Set-TextValue 'E41' '-3.18%'
Set-TextValue 'D42' '0.009544'
Set-TextValue 'E42' '22.26%'
Set-TextValue 'E43' '-5.54%'
Set-TextValue 'D44' '0.002110'
Set-TextValue 'E44' '0.72%'
Set-TextValue 'E45' '-7.77%'
Set-TextValue 'D46' '0.00006020'
Set-TextValue 'E46' '-6.93%'
Set-TextValue 'E47' '0.01%'
Set-TextValue 'E49' '-31.35%'
Set-TextValue 'D50' '0.00002100'
Set-TextValue 'E50' '0.01%'
Set-TextValue 'D51' '0.0002000'
Set-TextValue 'E51' '0.01%'
